# The target's "target" particle was renamed from "proton" to "p" in the
# data table (column G, rows 2-13), and the header row (row 1) was made bold.
# The previously-selected cell is also updated to match the author's final
# cursor position in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G = "target"; replace "proton" with "p" for every data row.
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 7).Value = "p"
}

# Bold the header row (row 1), keeping its existing center alignment.
$ws.Range("A1:K1").Font.Bold = $true

# Match the final selection left in the sheet.
$ws.Range("G18").Select()
